$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (matches source data
# where prices are strings even when they look like plain numbers), without leaving
# the cells visible style changed.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '34.711.66'
$ws.Range("E2").Value = '  +3.03%  '

$ws.Range("D3").Value = '1.789.89'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  -0.02%  '

Set-TextValue "D5" '223.17'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("E7").Value = '  +0.03%  '

Set-TextValue "D8" '32.58'
$ws.Range("E8").Value = '  +8.29%  '

$ws.Range("E9").Value = '  +1.10%  '

Set-TextValue "D10" '0.0688'
$ws.Range("E10").Value = '  +4.13%  '

$ws.Range("E11").Value = '  +1.56%  '

$ws.Range("D12").Value = '2.044.75'
$ws.Range("E12").Value = '  +0.86%  '

Set-TextValue "D13" '11.01'
$ws.Range("E13").Value = '  +9.83%  '

$ws.Range("D14").Value = '1.790.22'
$ws.Range("E14").Value = '  +0.97%  '

$ws.Range("D15").Value = '34.720.89'
$ws.Range("E15").Value = '  +3.15%  '

Set-TextValue "D16" '0.630'
$ws.Range("E16").Value = '  +1.26%  '

$ws.Range("E17").Value = '  +3.14%  '

Set-TextValue "D18" '68.52'
$ws.Range("E18").Value = '  +0.37%  '

Set-TextValue "D19" '253.19'
$ws.Range("E19").Value = '  +1.41%  '

$ws.Range("D20").Value = '0.0₃0787'
$ws.Range("E20").Value = '  +6.85%  '

$ws.Range("E21").Value = '  -0.24%  '

Set-TextValue "D22" '10.47'
$ws.Range("E22").Value = '  +2.22%  '

Set-TextValue "D24" '2.12'
$ws.Range("E24").Value = '  +0.14%  '

Set-TextValue "D25" '158.25'
$ws.Range("E25").Value = '  +0.25%  '

$ws.Range("E26").Value = '  -0.09%  '

Set-TextValue "D27" '7.05'
$ws.Range("E27").Value = '  +1.64%  '

$ws.Range("E28").Value = '  +0.17%  '

Set-TextValue "D30" '0.0516'
$ws.Range("E30").Value = '  +0.42%  '

$ws.Range("E31").Value = '  -1.00%  '

$ws.Range("E32").Value = '  -0.28%  '

Set-TextValue "D33" '3.57'
$ws.Range("E33").Value = '  +0.63%  '

Set-TextValue "D34" '1.88'
$ws.Range("E34").Value = '  +2.67%  '

$ws.Range("D35").Value = '1.432.75'
$ws.Range("E35").Value = '  -3.03%  '

Set-TextValue "D36" '1.05'
$ws.Range("E36").Value = '  -0.93%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D37" '0.0189'
$ws.Range("E37").Value = '  +2.60%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D38" '0.631'
$ws.Range("E38").Value = '  +0.75%  '

Set-TextValue "D39" '83.24'
$ws.Range("E39").Value = '  +0.43%  '

Set-TextValue "D40" '2.80'
$ws.Range("E40").Value = '  +4.33%  '

$ws.Range("E41").Value = '  +0.21%  '

Set-TextValue "D42" '0.902'
$ws.Range("E42").Value = '  +1.95%  '

$ws.Range("E43").Value = '  -0.58%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D44" '1.07'
$ws.Range("E44").Value = '  -0.57%  '

$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D45" '0.0504'
$ws.Range("E45").Value = '  -1.15%  '

$ws.Range("E46").Value = '  +4.37%  '

$ws.Range("D47").Value = '1.944.41'
$ws.Range("E47").Value = '  +1.34%  '

Set-TextValue "D48" '104.12'
$ws.Range("E48").Value = '  +7.12%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D49" '11.99'
$ws.Range("E49").Value = '  +3.21%  '

$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D50" '0.999'
$ws.Range("E50").Value = '  -0.06%  '

Set-TextValue "D51" '49.79'
$ws.Range("E51").Value = '  -1.96%  '
